$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "Row num: 0 Cell num: 2"
$ws.Range("C5").Value = "Row num: 1 Cell num: 2"
$ws.Range("C6").Value = "Row num: 2 Cell num: 2"
$ws.Range("C7").Value = "Row num: 3 Cell num: 2"

$ws.Range("C4:C7").Style = $ws.Range("B4").Style
